$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cases")

# Update the "management" column (G) and "rotation" column (F) for all data rows.
$ws.Range("G2:G6").Value = 'c("achille_rainfed_3N","achille_rainfed_3N")'
$ws.Range("F2:F6").Value = 'c("WHEAT.Ble_Dur_1", "Chickpea.Ghab2")'

# Widen column F to fit the new, longer text.
$ws.Columns.Item(6).ColumnWidth = 35.33

# Move the active selection (matches the saved file's cursor position).
$ws.Range("I7").Select()
